$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.369.95'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.812.49'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +5.46%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.56%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '343.90'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +3.54%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3809'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +3.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3501'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +4.71%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '49.16'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.51%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +3.52%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07727'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.23%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.02'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +9.64%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.624'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +5.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.263'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +4.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.808.92'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +5.60%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001119'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.78%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06726'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '86.01'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.66'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +7.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.560'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +7.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.29'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '27.384.29'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +5.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.466'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.670'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +7.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.07'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +14.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.478'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +12.14%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '153.47'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +2.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.009.18'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +5.48%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '136.35'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +5.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.321'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +6.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.047'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '13.87'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +7.39%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08736'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.50%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.709'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.627'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +4.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.6999'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +13.32%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2279'
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02417'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +5.49%  '
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.06503'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +4.72%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.981'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +5.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.302'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.94%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.77'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6511'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +10.50%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.003'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.40%  '
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +5.04%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.185'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +7.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '132.17'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +3.43%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07341'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '80.50'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.37%  '
